$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at L (pushing the existing "Ghi chú" column from L to M)
$ws.Columns("L").Insert()

# Fill in the new "Mã đơn hàng" column
$ws.Range("L1").Value = "Mã đơn hàng"
$ws.Range("L2").Value = "DH1"
$ws.Range("L3").Value = "DH2"

# Give the new column a custom width similar to column K's
$ws.Range("L1").ColumnWidth = 17.67

# Update the active selection to match the edited workbook
[void]$ws.Range("L6").Select()
